$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.255.31'
$ws.Range("E2").Value = '  +7.63%  '

# Row 3
$ws.Range("D3").Value = '3.343.03'
$ws.Range("E3").Value = '  +2.43%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '''412.32'
$ws.Range("E5").Value = '  +3.97%  '

# Row 6
$ws.Range("D6").Value = '''116.86'
$ws.Range("E6").Value = '  +7.09%  '

# Row 7
$ws.Range("D7").Value = '3.337.93'
$ws.Range("E7").Value = '  +2.31%  '

# Row 8
$ws.Range("E8").Value = '  -1.87%  '

# Row 9
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  -0.04%  '

# Row 10
$ws.Range("D10").Value = '''0.635'
$ws.Range("E10").Value = '  +1.06%  '

# Row 11
$ws.Range("E11").Value = '  +18.13%  '

# Row 12
$ws.Range("D12").Value = '''40.34'
$ws.Range("E12").Value = '  +2.48%  '

# Row 13
$ws.Range("E13").Value = '  -0.72%  '

# Row 14
$ws.Range("D14").Value = '3.865.62'
$ws.Range("E14").Value = '  +2.49%  '

# Row 15
$ws.Range("D15").Value = '''8.35'
$ws.Range("E15").Value = '  -1.06%  '

# Row 16
$ws.Range("D16").Value = '''19.28'
$ws.Range("E16").Value = '  +0.38%  '

# Row 17
$ws.Range("D17").Value = '3.352.64'
$ws.Range("E17").Value = '  +2.97%  '

# Row 18
$ws.Range("D18").Value = '61.006.85'
$ws.Range("E18").Value = '  +7.38%  '

# Row 19
$ws.Range("E19").Value = '  -2.32%  '

# Row 20
$ws.Range("D20").Value = '''10.89'
$ws.Range("E20").Value = '  +1.07%  '

# Row 21
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '''0.0000115'
$ws.Range("E21").Value = '  +5.42%  '

# Row 22
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").Value = '''3.38'
$ws.Range("E22").Value = '  +0.96%  '

# Row 23
$ws.Range("D23").Value = '''12.55'
$ws.Range("E23").Value = '  -3.70%  '

# Row 24
$ws.Range("D24").Value = '''297.79'
$ws.Range("E24").Value = '  +0.59%  '

# Row 25
$ws.Range("D25").Value = '''74.34'
$ws.Range("E25").Value = '  -0.15%  '

# Row 26
$ws.Range("E26").Value = '  -1.70%  '

# Row 27
$ws.Range("D27").Value = '''29.29'
$ws.Range("E27").Value = '  +3.83%  '

# Row 28
$ws.Range("D28").Value = '''7.82'
$ws.Range("E28").Value = '  +7.32%  '

# Row 29
$ws.Range("E29").Value = '  -2.69%  '

# Row 30
$ws.Range("D30").Value = '''0.172'
$ws.Range("E30").Value = '  +1.55%  '

# Row 31
$ws.Range("D31").Value = '''7.57'
$ws.Range("E31").Value = '  -1.53%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.115'
$ws.Range("E32").Value = '  +5.33%  '

# Row 33
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").Value = '''42.88'
$ws.Range("E33").Value = '  +7.21%  '

# Row 34
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D34").Value = '''2.54'
$ws.Range("E34").Value = '  +19.20%  '

# Row 35
$ws.Range("D35").Value = '''11.36'
$ws.Range("E35").Value = '  +0.70%  '

# Row 36
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  +0.01%  '

# Row 37
$ws.Range("D37").Value = '''0.0491'
$ws.Range("E37").Value = '  +0.13%  '

# Row 38
$ws.Range("D38").Value = '''52.42'
$ws.Range("E38").Value = '  +1.74%  '

# Row 39
$ws.Range("D39").Value = '''0.998'
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("E40").Value = '  +4.52%  '

# Row 41
$ws.Range("D41").Value = '''3.45'
$ws.Range("E41").Value = '  -1.01%  '

# Row 42
$ws.Range("D42").Value = '''134.71'
$ws.Range("E42").Value = '  -3.11%  '

# Row 43
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '''0.121'
$ws.Range("E43").Value = '  -0.83%  '

# Row 44
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '''0.290'
$ws.Range("E44").Value = '  +2.78%  '

# Row 45
$ws.Range("D45").Value = '''1.91'
$ws.Range("E45").Value = '  +0.00%  '

# Row 46
$ws.Range("D46").Value = '''3.92'
$ws.Range("E46").Value = '  -2.29%  '

# Row 47
$ws.Range("D47").Value = '''16.47'
$ws.Range("E47").Value = '  -3.94%  '

# Row 48
$ws.Range("E48").Value = '  +4.52%  '

# Row 49
$ws.Range("D49").Value = '''21.20'
$ws.Range("E49").Value = '  -4.89%  '

# Row 50
$ws.Range("D50").Value = '2.154.33'
$ws.Range("E50").Value = '  -0.22%  '

# Row 51
$ws.Range("D51").Value = '3.667.86'
$ws.Range("E51").Value = '  +2.56%  '
